$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row pairs whose match-data (everything except the running index column A,
# the Div column C and the Date column D) was swapped between the two rows.
$rowPairs = @(
    @(24, 25),
    @(26, 27),
    @(38, 39),
    @(41, 42),
    @(80, 81),
    @(98, 99),
    @(155, 156),
    @(182, 183)
)

# Columns that participate in the swap: B (id) plus E..AB (everything from
# HomeTeam through PL_AhUnder). Column numbers: B=2, E=5 .. AB=28.
$cols = @(2) + (5..28)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    foreach ($c in $cols) {
        $cell1 = $ws.Cells.Item($r1, $c)
        $cell2 = $ws.Cells.Item($r2, $c)
        $v1 = $cell1.Value2
        $v2 = $cell2.Value2
        $cell1.Value2 = $v2
        $cell2.Value2 = $v1
    }
}

# Standalone odds corrections on two not-yet-played fixtures.
$ws.Range("K202").Value2 = 3.2
$ws.Range("N202").Value2 = 3.2
$ws.Range("Q202").Value2 = 1.8
$ws.Range("R202").Value2 = 2

$ws.Range("T203").Value2 = 1.85
$ws.Range("U203").Value2 = 1.95
